# Updated cryptos list on Sun Aug 20 17:07:18 UTC 2023 with GitHub Actions
#
# This refreshes the "Price" (column D) and "Volume(1h)" (column E) figures
# for each coin row on the active worksheet, and fixes the ordering of the
# VeChain / Maker rows (37 / 38) which had been swapped.
#
# Note: several "Price" values (e.g. "26.336.46", "1.012") look like plain
# numbers to Excel's auto-detection, so they are written with a leading
# apostrophe to force them to stay text, exactly like the rest of the
# column. The apostrophe itself is not stored as part of the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.336.46"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.685.98"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.79%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'218.22"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.5403"
$ws.Range("E6").Value = "  +2.70%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.75%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2731"
$ws.Range("E8").Value = "  +1.00%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06442"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'21.88"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +2.72%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "'1.694.80"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'4.525"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.5797"
$ws.Range("E14").Value = "  -1.57%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "'0.000008345"
$ws.Range("E15").Value = "  -3.31%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'64.95"
$ws.Range("E16").Value = "  -0.04%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'26.413.09"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'4.920"
$ws.Range("E18").Value = "  -1.79%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.75%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'10.95"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'191.06"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "'6.235"
$ws.Range("E22").Value = "  -0.62%  "

# Row 23 - BinanceUSD
$ws.Range("D23").Value = "'1.012"
$ws.Range("E23").Value = "  +0.69%  "

# Row 24 - Monero
$ws.Range("D24").Value = "'149.54"
$ws.Range("E24").Value = "  +2.80%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1317"
$ws.Range("E25").Value = "  +6.05%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.841"
$ws.Range("E26").Value = "  +1.68%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'15.72"
$ws.Range("E27").Value = "  -1.37%  "

# Row 28 - Hedera
$ws.Range("D28").Value = "'0.06329"
$ws.Range("E28").Value = "  -7.50%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +4.70%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.325"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "'3.578"
$ws.Range("E31").Value = "  -1.31%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'3.561"
$ws.Range("E32").Value = "  -0.39%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "'1.683"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +0.31%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "'0.6124"
$ws.Range("E35").Value = "  -2.02%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "'2.414"
$ws.Range("E36").Value = "  +1.34%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = "  -0.12%  "

# Row 38 - FraxShare
$ws.Range("D38").Value = "'6.262"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39 - was VeChain, now Maker (rows 39/40 swapped)
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.108.76"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40 - was Maker, now VeChain (rows 39/40 swapped)
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01624"
$ws.Range("E40").Value = "  -0.01%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.8811"
$ws.Range("E41").Value = "  +0.50%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.10%  "

# Row 43 - Quant
$ws.Range("D43").Value = "'101.78"
$ws.Range("E43").Value = "  +0.71%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "'1.837.48"
$ws.Range("E44").Value = "  -0.41%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  -0.17%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'57.35"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47 - Frax
$ws.Range("D47").Value = "'1.011"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'8.187"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.05266"
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "'0.4308"
$ws.Range("E50").Value = "  +0.27%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "'6.023"
$ws.Range("E51").Value = "  -0.63%  "
